$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data stores every Coin/Link/Price/Volume cell as plain text (t="inlineStr"),
# even when a price like "592.93" or "1.00" looks numeric. A bare .Value assignment
# lets Excel auto-detect such strings as numbers (dropping the text type, and for
# trailing-zero prices like "1.00"/"1.30"/"1.70" even rewriting the digits). Prefix
# those cells with an apostrophe, just like a user forcing text entry in the UI, so
# the literal string is preserved.

$ws.Range("D2").Value = "68.714.78"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "2.515.80"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'592.93"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Value = "'174.27"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "2.515.43"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "'4.99"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'0.336"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "2.979.43"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "'25.82"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "68.584.19"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "'0.0000172"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "2.516.57"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").Value = "'362.87"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").Value = "'7.52"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'10.87"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'70.29"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "'4.16"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'8.86"
$ws.Range("E26").Value = "  -3.75%  "
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "'1.65"
$ws.Range("E27").Value = "  -7.55%  "
$ws.Range("D28").Value = "2.646.94"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'508.54"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "0.0₃0872"
$ws.Range("E31").Value = "  -4.35%  "
$ws.Range("D32").Value = "'7.72"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'161.57"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "'18.56"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'18.64"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'1.30"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").Value = "'1.70"
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").Value = "'0.322"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "'4.74"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "'2.34"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").Value = "'150.32"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'0.511"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "'0.0736"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "0.0₆0249"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'1.55"
$ws.Range("E51").Value = "  -2.10%  "
